$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.609.22"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.432.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.56%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.28%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.28%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.73%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("E9").Value = "  +2.77%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.98%  "

# Row 12
$ws.Range("E12").Value = "  +1.05%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.93%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.812.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.44%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.408.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.79%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.838"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.60%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.557.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.74%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.50%  "

# Row 20
$ws.Range("E20").Value = "  +1.59%  "

# Row 21
$ws.Range("E21").Value = "  +2.30%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.85%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.55%  "

# Row 25
$ws.Range("E25").Value = "  +1.77%  "

# Row 26
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.61%  "

# Row 28
$ws.Range("E28").Value = "  -4.17%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.67"
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.73%  "

# Row 31
$ws.Range("B31").Value = "OKB"
$ws.Range("C31").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.33%  "

# Row 32
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +16.47%  "

# Row 33
$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.78%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.81%  "

# Row 35
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0761"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.75%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.77%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.34%  "

# Row 39
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.03%  "

# Row 40
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.37%  "

# Row 41
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.109"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.83%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.54%  "

# Row 43
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.08%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0291"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.61%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.950.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.41%  "

# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.19%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.18%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.36%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.48%  "

# Row 50
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.02%  "

# Row 51
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.76%  "
